$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell J1 = "Q8", copying style (bold/border) from I1
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J1").Value = "Q8"

# Update numeric data (simulated rt_data bugfix)
# Row 2
$ws.Range("B2").Value = 2.328981359987906
$ws.Range("C2").Value = -0.6012016438619422
$ws.Range("D2").Value = -0.837771994326261
$ws.Range("E2").Value = -0.5121775623451877
$ws.Range("F2").Value = 0.05154791033427741
$ws.Range("G2").Value = -0.2346076846697836
$ws.Range("H2").Value = -0.1127621191487843

# Row 3
$ws.Range("B3").Value = -0.9720654395460997
$ws.Range("C3").Value = -1.208635790010419
$ws.Range("D3").Value = -0.8830413580293452
$ws.Range("E3").Value = -0.3193158853498801
$ws.Range("F3").Value = -0.605471480353941
$ws.Range("G3").Value = -0.4836259148329418

# Row 4
$ws.Range("B4").Value = -0.7302303441561699
$ws.Range("C4").Value = -0.4046359121750967
$ws.Range("D4").Value = 0.1590895605043684
$ws.Range("E4").Value = -0.1270660344996926
$ws.Range("F4").Value = -0.005220468978693271
$ws.Range("G4").Value = -0.4389000132169812
$ws.Range("H4").Value = 1.446976225670126
$ws.Range("I4").Value = -0.2565828398263235
$ws.Range("J4").Value = -0.3121643713505491

# Row 5
$ws.Range("B5").Value = -0.2663886325293134
$ws.Range("C5").Value = 0.2973368401501517
$ws.Range("D5").Value = 0.0111812451460907
$ws.Range("E5").Value = 0.13302681066709
$ws.Range("F5").Value = -0.3006527335711979
$ws.Range("G5").Value = 1.585223505315909
$ws.Range("H5").Value = -0.1183355601805403
$ws.Range("I5").Value = -0.1739170917047659

# Row 6
$ws.Range("B6").Value = 0.2463078550078095
$ws.Range("C6").Value = -0.03984773999625146
$ws.Range("D6").Value = 0.08199782552474782
$ws.Range("E6").Value = -0.3516817187135401
$ws.Range("F6").Value = 1.534194520173567
$ws.Range("G6").Value = -0.1693645453228824
$ws.Range("H6").Value = -0.224946076847108

# Row 7
$ws.Range("B7").Value = -0.07311312868471159
$ws.Range("C7").Value = 0.0487324368362877
$ws.Range("D7").Value = -0.3849471074020002
$ws.Range("E7").Value = 1.500929131485107
$ws.Range("F7").Value = -0.2026299340113425
$ws.Range("G7").Value = -0.2582114655355682

# Row 8
$ws.Range("B8").Value = 0.1900396483956045
$ws.Range("C8").Value = -0.2436398958426834
$ws.Range("D8").Value = 1.642236343044424
$ws.Range("E8").Value = -0.06132272245202575
$ws.Range("F8").Value = -0.1169042539762514
$ws.Range("G8").Value = 0.6700456085978389
$ws.Range("H8").Value = -1.519651802199725
$ws.Range("I8").Value = -0.6816607491277794

# Row 9
$ws.Range("B9").Value = -0.3722755437654697
$ws.Range("C9").Value = 1.513600695121637
$ws.Range("D9").Value = -0.1899583703748121
$ws.Range("E9").Value = -0.2455399018990377
$ws.Range("F9").Value = 0.5414099606750526
$ws.Range("G9").Value = -1.648287450122512
$ws.Range("H9").Value = -0.8102963970505657

# Row 10
$ws.Range("B10").Value = 1.52334500150552
$ws.Range("C10").Value = -0.1802140639909298
$ws.Range("D10").Value = -0.2357955955151554
$ws.Range("E10").Value = 0.5511542670589349
$ws.Range("F10").Value = -1.638543143738629
$ws.Range("G10").Value = -0.8005520906666834

# Row 11
$ws.Range("B11").Value = -0.4633150438766606
$ws.Range("C11").Value = -0.5188965754008863
$ws.Range("D11").Value = 0.268053287173204
$ws.Range("E11").Value = -1.92164412362436
$ws.Range("F11").Value = -1.083653070552414

# Row 12
$ws.Range("B12").Value = -0.2040875226113216
$ws.Range("C12").Value = 0.5828623399627687
$ws.Range("D12").Value = -1.606835070834796
$ws.Range("E12").Value = -0.7688440177628496

# Row 13
$ws.Range("B13").Value = 0.6123731026384148
$ws.Range("C13").Value = -1.577324308159149
$ws.Range("D13").Value = -0.7393332550872035

# Row 14
$ws.Range("B14").Value = -1.714346676517737
$ws.Range("C14").Value = -0.8763556234457911

# Row 15
$ws.Range("B15").Value = -0.7422840489468245
